$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Move to location (7, 5) and remove the screws."
$ws.Range("B2").Value = 22.767147
$ws.Range("C2").Value = 3837
$ws.Range("D2").Value = "'0.00804"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "a2c581be-a93b-4358-aba4-e81a987ed28b"

# Row 3
$ws.Range("A3").Value = "Move Robot48 to location (6, 5) and remove the liquid spill."
$ws.Range("B3").Value = 42.24573
$ws.Range("C3").Value = 5950
$ws.Range("D3").Value = "'0.01353"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "040e981e-b63f-4ccd-b75f-44796766e9c8"

# Row 4
$ws.Range("A4").Value = "Move Robot35 to location (3, 3) and remove the large debris."
$ws.Range("B4").Value = 28.530881
$ws.Range("C4").Value = 4464
$ws.Range("D4").Value = "'0.00888"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "44bb4645-a38c-4a5c-be61-f4a96e197333"

# Row 5
$ws.Range("A5").Value = "Move Robot15 to location (3, 4) and remove the dust."
$ws.Range("B5").Value = 29.356178
$ws.Range("C5").Value = 4448
$ws.Range("D5").Value = "'0.00906"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "0483bbe5-8a4e-4562-9c79-c6b42f8bf898"

# Row 6
$ws.Range("A6").Value = "Move Robot29 to location (9, 10) and remove the grass."
$ws.Range("B6").Value = 30.561105
$ws.Range("C6").Value = 4559
$ws.Range("D6").Value = "'0.0096"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "56d5ed3c-117f-436c-89cd-9977af779654"

# Row 7
$ws.Range("A7").Value = "Move Robot31 to location (8, 12) and remove the small debris."
$ws.Range("B7").Value = 28.585955
$ws.Range("C7").Value = 4484
$ws.Range("D7").Value = "'0.00903"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "bf9becbe-829f-4d58-9eac-ad266fdf7e47"

# Row 8
$ws.Range("A8").Value = "Move Robot13 to location (7, 5) and remove the vehicle."
$ws.Range("B8").Value = 29.192122
$ws.Range("C8").Value = 4557
$ws.Range("D8").Value = "'0.00966"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "5b8e0628-634e-4827-985e-d92b7e4d8cdf"

# Row 9
$ws.Range("A9").Value = "Move Robot50 to location (5, 12) and remove the construction materials."
$ws.Range("B9").Value = 85.219786
$ws.Range("C9").Value = 19963
$ws.Range("D9").Value = "'0.027"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "624f9048-ab6b-4977-95ba-23ba7f9f3184"

# Row 10
$ws.Range("A10").Value = "Move Robot9 to location (11, 2) and remove the tree branches."
$ws.Range("B10").Value = 28.737158
$ws.Range("C10").Value = 4436
$ws.Range("D10").Value = "'0.00933"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "b0d4f7b4-f07d-45b3-9006-2bc464ed0909"

# Row 11
$ws.Range("A11").Value = "Move Robot40 to location (10, 3) and remove the screws."
$ws.Range("B11").Value = 116.261215
$ws.Range("C11").Value = 36166
$ws.Range("D11").Value = "'0.04173"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "505fef53-6169-4d85-acf9-11a98bb67cfd"
